# Update NATMI LR-pairs output (C3-Lrp1) with recomputed values based on new TPM data.
# The underlying TPM matrix used to derive Ligand/Receptor expression values for the
# "ECs" cluster changed, which cascades through the detection rate, average/total
# expression, derived-specificity, and edge-weight columns for every row touching ECs
# either as a sending (source) or target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.04155
$ws.Range("H2").Value = 0.12465
$ws.Range("I2").Value = 0.0001466168179836329
$ws.Range("J2").Value = 0.0001466168179836329
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 0.1436078246
$ws.Range("R2").Value = 1.2924704214
$ws.Range("S2").Value = 0.000001442934664193658
$ws.Range("T2").Value = 0.000001442934664193657

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.04155
$ws.Range("H3").Value = 0.12465
$ws.Range("I3").Value = 0.0001466168179836329
$ws.Range("J3").Value = 0.0001466168179836329
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 12.5314861771
$ws.Range("R3").Value = 112.7833755939
$ws.Range("S3").Value = 0.0001259131655894553
$ws.Range("T3").Value = 0.0001259131655894552

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.04155
$ws.Range("H4").Value = 0.12465
$ws.Range("I4").Value = 0.0001466168179836329
$ws.Range("J4").Value = 0.0001466168179836329
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 1.91691962365
$ws.Range("R4").Value = 17.25227661285
$ws.Range("S4").Value = 0.00001926071772998395
$ws.Range("T4").Value = 0.00001926071772998395

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.9992428949822291
$ws.Range("J5").Value = 0.9992428949822291
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 978.7355936985567
$ws.Range("R5").Value = 8808.620343287012
$ws.Range("S5").Value = 0.009834084731534937
$ws.Range("T5").Value = 0.009834084731534935

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.9992428949822291
$ws.Range("J6").Value = 0.9992428949822291
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.8581405450637278
$ws.Range("T6").Value = 0.8581405450637277

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.9992428949822291
$ws.Range("J7").Value = 0.9992428949822291
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.1312682651869665
$ws.Range("T7").Value = 0.1312682651869664

# Row 8: MuSCs -> ECs
$ws.Range("I8").Value = 0.0006104881997874136
$ws.Range("J8").Value = 0.0006104881997874135
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 0.5979592486124443
$ws.Range("R8").Value = 5.381633237511999
$ws.Range("S8").Value = 0.000006008141478372407
$ws.Range("T8").Value = 0.000006008141478372405

# Row 9: MuSCs -> FAPs
$ws.Range("I9").Value = 0.0006104881997874136
$ws.Range("J9").Value = 0.0006104881997874135
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.0005242816127602909
$ws.Range("T9").Value = 0.0005242816127602907

# Row 10: MuSCs -> MuSCs
$ws.Range("I10").Value = 0.0006104881997874136
$ws.Range("J10").Value = 0.0006104881997874135
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.00008019844554875035
$ws.Range("T10").Value = 0.00008019844554875033
